$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 183, pushing existing rows 183:193 down to 184:194
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with the new record
$ws.Cells.Item(183, 1).Value = 5
$ws.Cells.Item(183, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(183, 3).Value = "Maule"
$ws.Cells.Item(183, 4).Value = 44746
$ws.Cells.Item(183, 4).NumberFormat = $ws.Cells.Item(184, 4).NumberFormat
$ws.Cells.Item(183, 5).Value = 7
$ws.Cells.Item(183, 6).Value = 100112017
$ws.Cells.Item(183, 7).Value = "Apio"
$ws.Cells.Item(183, 8).Value = "Americana (o)"
$ws.Cells.Item(183, 9).Value = "Primera"
$ws.Cells.Item(183, 10).Value = 600
$ws.Cells.Item(183, 11).Value = 6500
$ws.Cells.Item(183, 12).Value = 6500
$ws.Cells.Item(183, 13).Value = 6500
$ws.Cells.Item(183, 14).Value = "`$/docena de matas"
$ws.Cells.Item(183, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(183, 16).Value = 1083
$ws.Cells.Item(183, 17).Value = 6
$ws.Cells.Item(183, 18).Value = "Hortaliza"
